# Regenerate merged AHB files
# 1. Rename header labels: *_old -> *_FV2404, *_new -> *_FV2410
# 2. Add a table (ListObject) over the used range A1:U65
# 3. Freeze the header row (pane split) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2404 = "_FV2404"
$fv2410 = "_FV2410"

$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

# Columns A-J -> *_old => *_FV2404
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value2 = $headers[$i] + $fv2404
}

# Column K = "diff" stays the same

# Columns L-U -> *_new => *_FV2410
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 12)
    $cell.Value2 = $headers[$i] + $fv2410
}

# Add table over the full used range (header row + all data rows)
$rng = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $rng, [Type]::Missing, 1)
$tbl.Name = "Table1"

# Freeze header row (split pane after row 1)
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
